$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows (28..102) down to (29..103)
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly price record
$ws.Cells.Item(28, 1).Value  = 7
$ws.Cells.Item(28, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(28, 3).Value  = "Ñuble"
$ws.Cells.Item(28, 4).Value  = 44868
$ws.Cells.Item(28, 5).Value  = 16
$ws.Cells.Item(28, 6).Value  = 100112031
$ws.Cells.Item(28, 7).Value  = "Poroto verde"
$ws.Cells.Item(28, 8).Value  = "Magnum"
$ws.Cells.Item(28, 9).Value  = "Primera"
$ws.Cells.Item(28, 10).Value = 60
$ws.Cells.Item(28, 11).Value = 35000
$ws.Cells.Item(28, 12).Value = 36000
$ws.Cells.Item(28, 13).Value = 35500
$ws.Cells.Item(28, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value = 1420
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"
